$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, more detailed / realistic Indonesian app-review comments replacing
# the old sample comments, plus several brand-new rows.
$comments = @(
  "Aplikasinya sangat bagus dan membantu pekerjaan saya!",
  "Lemot banget, sering keluar sendiri pas lagi dipakai.",
  "Desainnya keren sih, tapi kadang masih suka error.",
  "Suka banget sama update terbarunya, jadi lebih cepat.",
  "CS-nya parah, gak ada yang bales email komplain saya.",
  "Lumayan lah, tapi iklannya kebanyakan jadi ganggu.",
  "Aplikasi terbaik untuk produktivitas, sangat rekomen!",
  "Gak bisa login setelah ganti password. Tolong diperbaiki.",
  "Fiturnya berguna, tapi tampilannya agak membingungkan.",
  "Kasih bintang lima! Berfungsi sesuai deskripsi."
)

for ($i = 0; $i -lt $comments.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $comments[$i]
}

# Uniform look for every comment row: Arial 10, vertically centered (no
# wrap anymore), boxed in with a light-grey medium border.
$dataRange = $ws.Range("A2:A11")
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $false
$dataRange.Font.Name = "Arial"
$dataRange.Font.Size = 10
$dataRange.Borders.Weight = -4138
$dataRange.Borders.Color = 0xCCCCCC

# Tighter fixed row heights across the whole table.
$ws.Range("A1:A11").Rows.RowHeight = 15

# Widen the column so the longer review text fits.
$ws.Columns.Item(1).ColumnWidth = 47.36328125

# Match the saved selection/active cell.
[void]$ws.Range("A2:A11").Select()
